$p = $ppt.ActivePresentation

# --- Change 1: slide 3, "TextBox 1" shape ---
# Paragraph (lvl=1) originally held two runs: a lone ko-KR " " run followed by
# the en-US run with "However, current methods...". Remove the leading
# ko-KR " " run entirely so only the en-US run (and its text) remains.
$s3 = $p.Slides.Item(3)
$shp3 = $s3.Shapes.Item(2)
$tr3 = $shp3.TextFrame.TextRange
$para3 = $tr3.Paragraphs(7, 1)
$lead3 = $para3.Characters(1, 1)
$lead3.Text = ""

# --- Change 2: slide 5, "TextBox 37" shape ---
# Last paragraph's text had a stray leading space before "Is the cheek...".
# Strip that leading space while leaving the run's formatting untouched.
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(6)
$tr5 = $shp5.TextFrame.TextRange
$para5 = $tr5.Paragraphs(10, 1)
$lead5 = $para5.Characters(1, 1)
$lead5.Text = ""
